$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: snapshot old F:V values for all rows being relocated (reordering within date groups) ---
$snap = @{}
$snap[76] = $ws.Range("F76:V76").Value()
$snap[78] = $ws.Range("F78:V78").Value()
$snap[85] = $ws.Range("F85:V85").Value()
$snap[86] = $ws.Range("F86:V86").Value()
$snap[87] = $ws.Range("F87:V87").Value()
$snap[88] = $ws.Range("F88:V88").Value()
$snap[89] = $ws.Range("F89:V89").Value()
$snap[96] = $ws.Range("F96:V96").Value()
$snap[97] = $ws.Range("F97:V97").Value()
$snap[98] = $ws.Range("F98:V98").Value()
$snap[99] = $ws.Range("F99:V99").Value()
$snap[106] = $ws.Range("F106:V106").Value()
$snap[107] = $ws.Range("F107:V107").Value()
$snap[108] = $ws.Range("F108:V108").Value()
$snap[109] = $ws.Range("F109:V109").Value()
$snap[110] = $ws.Range("F110:V110").Value()
$snap[111] = $ws.Range("F111:V111").Value()
$snap[113] = $ws.Range("F113:V113").Value()
$snap[114] = $ws.Range("F114:V114").Value()
$snap[115] = $ws.Range("F115:V115").Value()
$snap[120] = $ws.Range("F120:V120").Value()
$snap[121] = $ws.Range("F121:V121").Value()
$snap[122] = $ws.Range("F122:V122").Value()
$snap[123] = $ws.Range("F123:V123").Value()
$snap[126] = $ws.Range("F126:V126").Value()
$snap[127] = $ws.Range("F127:V127").Value()
$snap[128] = $ws.Range("F128:V128").Value()
$snap[129] = $ws.Range("F129:V129").Value()
$snap[130] = $ws.Range("F130:V130").Value()

# --- Step 2: write snapshots into their destination rows ---
$ws.Range("F78:V78").Value = $snap[76]
$ws.Range("F76:V76").Value = $snap[78]
$ws.Range("F86:V86").Value = $snap[85]
$ws.Range("F85:V85").Value = $snap[86]
$ws.Range("F88:V88").Value = $snap[87]
$ws.Range("F89:V89").Value = $snap[88]
$ws.Range("F87:V87").Value = $snap[89]
$ws.Range("F97:V97").Value = $snap[96]
$ws.Range("F98:V98").Value = $snap[97]
$ws.Range("F99:V99").Value = $snap[98]
$ws.Range("F96:V96").Value = $snap[99]
$ws.Range("F107:V107").Value = $snap[106]
$ws.Range("F109:V109").Value = $snap[107]
$ws.Range("F110:V110").Value = $snap[108]
$ws.Range("F108:V108").Value = $snap[109]
$ws.Range("F106:V106").Value = $snap[110]
$ws.Range("F114:V114").Value = $snap[111]
$ws.Range("F115:V115").Value = $snap[113]
$ws.Range("F111:V111").Value = $snap[114]
$ws.Range("F113:V113").Value = $snap[115]
$ws.Range("F123:V123").Value = $snap[120]
$ws.Range("F120:V120").Value = $snap[121]
$ws.Range("F121:V121").Value = $snap[122]
$ws.Range("F122:V122").Value = $snap[123]
$ws.Range("F127:V127").Value = $snap[126]
$ws.Range("F126:V126").Value = $snap[127]
$ws.Range("F130:V130").Value = $snap[128]
$ws.Range("F128:V128").Value = $snap[129]
$ws.Range("F129:V129").Value = $snap[130]

# --- Step 3: append new match rows 145-151 (copy formats from row 144, then set values) ---
$ws.Range("A144:V144").Copy() | Out-Null
$ws.Range("A145:V151").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats) | Out-Null
$excel.CutCopyMode = $false

function Set-RowValues {
    param($wsArg, $rowNum, $values)
    $n = $values.Count
    $arr = New-Object 'object[,]' 1,$n
    for ($i = 0; $i -lt $n; $i++) {
        $arr[0,$i] = $values[$i]
    }
    $endCol = [char](64 + $n)
    $rng = "A" + $rowNum + ":" + $endCol + $rowNum
    $wsArg.Range($rng).Value = $arr
}

$vals145 = @(144, "greece", "super-league-2", "2023-2024", 45298.54166666666, "Iraklis 1908", 1, "AEL Larissa", 1, 3.07, "06/01/2024 01:12", 3.91, "07/01/2024 12:47", 2.7, "06/01/2024 01:12", 3.16, "07/01/2024 12:56", 2.44, "06/01/2024 01:12", 2, "07/01/2024 12:56", "https://www.betexplorer.com/football/greece/super-league-2/iraklis-fc-ael-larissa/boZpRIUe/")
Set-RowValues $ws 145 $vals145

$vals146 = @(145, "greece", "super-league-2", "2023-2024", 45298.5625, "Giouchtas", 1, "Kalamata", 1, 4.53, "06/01/2024 01:42", 4.57, "07/01/2024 13:17", 2.84, "06/01/2024 01:42", 3.12, "07/01/2024 13:17", 1.87, "06/01/2024 01:42", 1.88, "07/01/2024 13:17", "https://www.betexplorer.com/football/greece/super-league-2/giouchtas-kalamata/ne7jkmB9/")
Set-RowValues $ws 146 $vals146

$vals147 = @(146, "greece", "super-league-2", "2023-2024", 45298.58333333334, "Ionikos", 4, "Ilioupoli", 1, 1.74, "06/01/2024 02:12", 1.63, "07/01/2024 13:22", 3.22, "06/01/2024 02:12", 3.54, "07/01/2024 13:22", 4.55, "06/01/2024 02:12", 5.72, "07/01/2024 13:22", "https://www.betexplorer.com/football/greece/super-league-2/ionikos-ilioupoli/bZQeT5lj/")
Set-RowValues $ws 147 $vals147

$vals148 = @(147, "greece", "super-league-2", "2023-2024", 45298.58333333334, "PAE Egaleo", 4, "Diagoras", 1, 2.16, "06/01/2024 02:12", 2.11, "07/01/2024 13:24", 2.83, "06/01/2024 02:12", 2.95, "07/01/2024 13:24", 3.43, "06/01/2024 02:12", 3.87, "07/01/2024 13:24", "https://www.betexplorer.com/football/greece/super-league-2/pae-egaleo-diagoras-fc/Qo6fl7QF/")
Set-RowValues $ws 148 $vals148

$vals149 = @(148, "greece", "super-league-2", "2023-2024", 45298.58333333334, "Tilikratis L.", 0, "PAE Chania", 3, 8.61, "06/01/2024 02:12", 9.17, "07/01/2024 13:51", 4.14, "06/01/2024 02:12", 3.93, "07/01/2024 13:51", 1.34, "06/01/2024 02:12", 1.42, "07/01/2024 13:51", "https://www.betexplorer.com/football/greece/super-league-2/tilikratis-lefkada-pae-chania/hY7njTe3/")
Set-RowValues $ws 149 $vals149

$vals150 = @(149, "greece", "super-league-2", "2023-2024", 45298.58333333334, "Aiolikos", 1, "Karditsa", 1, 2.42, "06/01/2024 02:12", 2.01, "07/01/2024 13:59", 2.72, "06/01/2024 02:12", 2.99, "07/01/2024 13:59", 3.08, "06/01/2024 02:12", 4.16, "07/01/2024 13:59", "https://www.betexplorer.com/football/greece/super-league-2/aiolikos-fc-karditsa/4OZtSbFk/")
Set-RowValues $ws 150 $vals150

$vals151 = @(150, "greece", "super-league-2", "2023-2024", 45298.58333333334, "Makedonikos", 0, "Kozani FC", 0, 2.12, "06/01/2024 02:12", 2.39, "07/01/2024 12:03", 2.79, "06/01/2024 02:12", 2.81, "07/01/2024 12:03", 3.63, "06/01/2024 02:12", 3.36, "07/01/2024 12:03", "https://www.betexplorer.com/football/greece/super-league-2/makedonikos-neapolis-kozani-fc/pvgAXKiS/")
Set-RowValues $ws 151 $vals151
